$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 243569
$ws.Range("C5").Value = 25
$ws.Range("D5").Value = 42.63934426229508
$ws.Range("F5").Value = "2025-04-28 07:25:00"
$ws.Range("G5").Value = "2025-04-28 07:25:00"
$ws.Range("H5").Value = "2025-04-28 08:07:38"
$ws.Range("I5").Value = 2601
$ws.Range("K5").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;R12 ;R9"
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 76
$ws.Range("N5").Value = 3

# Row 6
$ws.Range("C6").Value = 42
$ws.Range("E6").Value = "2025-04-28 08:07:38"
$ws.Range("F6").Value = "2025-04-28 08:49:38"
$ws.Range("G6").Value = "2025-04-28 08:49:38"
$ws.Range("H6").Value = "2025-04-28 11:40:37"

# Row 7
$ws.Range("A7").Value = 251053
$ws.Range("C7").Value = 25
$ws.Range("D7").Value = 187.7377049180328
$ws.Range("E7").Value = "2025-04-28 11:40:37"
$ws.Range("F7").Value = "2025-04-28 12:05:37"
$ws.Range("G7").Value = "2025-04-28 12:05:37"
$ws.Range("H7").Value = "2025-04-29 07:13:21"
$ws.Range("I7").Value = 11452
$ws.Range("K7").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L7").Value = 4
$ws.Range("M7").Value = 70

# Row 8
$ws.Range("A8").Value = 251301
$ws.Range("C8").Value = 27
$ws.Range("D8").Value = 37.62295081967213
$ws.Range("E8").Value = "2025-04-29 07:13:21"
$ws.Range("F8").Value = "2025-04-29 07:40:21"
$ws.Range("G8").Value = "2025-04-29 07:40:21"
$ws.Range("H8").Value = "2025-04-29 08:17:59"
$ws.Range("I8").Value = 2295
$ws.Range("L8").Value = 5
$ws.Range("N8").ClearContents()

# Row 9
$ws.Range("A9").Value = 251182
$ws.Range("C9").Value = 34
$ws.Range("D9").Value = 302.3380281690141
$ws.Range("F9").Value = "2025-04-28 07:34:00"
$ws.Range("G9").Value = "2025-04-28 07:34:00"
$ws.Range("H9").Value = "2025-04-28 12:36:20"
$ws.Range("I9").Value = 21466
$ws.Range("L9").Value = 4

# Row 10
$ws.Range("A10").Value = 251300
$ws.Range("D10").Value = 54.12676056338028
$ws.Range("E10").Value = "2025-04-28 12:36:20"
$ws.Range("F10").Value = "2025-04-28 12:53:20"
$ws.Range("G10").Value = "2025-04-28 12:53:20"
$ws.Range("H10").Value = "2025-04-28 13:47:27"
$ws.Range("I10").Value = 3843
$ws.Range("L10").Value = 5

# Row 11
$ws.Range("C11").Value = 17
$ws.Range("E11").Value = "2025-04-28 13:47:27"
$ws.Range("F11").Value = "2025-04-28 14:04:27"
$ws.Range("G11").Value = "2025-04-28 14:04:27"
$ws.Range("H11").Value = "2025-04-29 08:04:33"

# Row 17
$ws.Range("A17").Value = 251308
$ws.Range("D17").Value = 81.9672131147541
$ws.Range("H17").Value = "2025-04-24 08:41:58"
$ws.Range("I17").Value = 5000

# Row 18
$ws.Range("A18").Value = 251168
$ws.Range("D18").Value = 142.3114754098361
$ws.Range("E18").Value = "2025-04-24 08:41:58"
$ws.Range("F18").Value = "2025-04-24 09:01:58"
$ws.Range("G18").Value = "2025-04-24 09:01:58"
$ws.Range("H18").Value = "2025-04-24 11:24:16"
$ws.Range("I18").Value = 8681

# Row 19
$ws.Range("A19").Value = 251167
$ws.Range("D19").Value = 173.655737704918
$ws.Range("E19").Value = "2025-04-24 11:24:16"
$ws.Range("F19").Value = "2025-04-24 11:44:16"
$ws.Range("G19").Value = "2025-04-24 11:44:16"
$ws.Range("I19").Value = 10593

# Row 20
$ws.Range("A20").Value = 250946
$ws.Range("C20").Value = 20
$ws.Range("D20").Value = 181.4098360655738
$ws.Range("F20").Value = "2025-04-24 14:57:56"
$ws.Range("G20").Value = "2025-04-24 14:57:56"
$ws.Range("H20").Value = "2025-04-25 09:59:20"
$ws.Range("I20").Value = 11066
$ws.Range("L20").Value = 4

# Row 21
$ws.Range("A21").Value = 250448
$ws.Range("D21").Value = 247.2622950819672
$ws.Range("E21").Value = "2025-04-25 09:59:20"
$ws.Range("F21").Value = "2025-04-25 10:29:20"
$ws.Range("G21").Value = "2025-04-25 10:29:20"
$ws.Range("H21").Value = "2025-04-25 14:36:36"
$ws.Range("I21").Value = 15083
$ws.Range("L21").Value = 2

# Row 22
$ws.Range("E22").Value = "2025-04-25 14:36:36"
$ws.Range("F22").Value = "2025-04-28 07:01:36"
$ws.Range("G22").Value = "2025-04-28 07:01:36"
$ws.Range("H22").Value = "2025-04-28 12:29:35"

# Row 23
$ws.Range("E23").Value = "2025-04-28 12:29:35"
$ws.Range("F23").Value = "2025-04-28 13:04:35"
$ws.Range("G23").Value = "2025-04-28 13:04:35"
$ws.Range("H23").Value = "2025-04-29 07:24:22"

# Row 29
$ws.Range("A29").Value = 251550
$ws.Range("C29").Value = 37
$ws.Range("D29").Value = 727.5714285714286
$ws.Range("F29").Value = "2025-04-28 07:37:00"
$ws.Range("G29").Value = "2025-04-28 07:37:00"
$ws.Range("H29").Value = "2025-04-29 11:44:34"
$ws.Range("I29").Value = 35651
$ws.Range("L29").Value = 3
$ws.Range("N29").Value = 4

# Row 30
$ws.Range("A30").Value = 251334
$ws.Range("D30").Value = 377.3877551020408
$ws.Range("E30").Value = "2025-04-29 11:44:34"
$ws.Range("F30").Value = "2025-04-29 12:24:34"
$ws.Range("G30").Value = "2025-04-29 12:24:34"
$ws.Range("H30").Value = "2025-04-30 10:41:57"
$ws.Range("I30").Value = 18492
$ws.Range("N30").ClearContents()
